$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 10773895
$ws.Range("I19").Value = 7224762
$ws.Range("K19").Value = 7224762
$ws.Range("M19").Value = -7224587

# Row 57: Quit Your Jib-jab | Gold Needle
$ws.Range("H57").Value = 28762.223
$ws.Range("J57").Value = 28762.223
$ws.Range("L57").Value = 86286.66900000001
$ws.Range("N57").Value = -87284.66900000001

# Row 70: Consecrating Congregation | Holy Water
$ws.Range("H70").Value = 3705506.8
$ws.Range("I70").Value = 1034.3334
$ws.Range("J70").Value = 5557743
$ws.Range("K70").Value = 3103.0002
$ws.Range("L70").Value = 16673229
$ws.Range("M70").Value = -2833.0002
$ws.Range("N70").Value = -16673769

# Row 73: Curbing the Contagion (L) | Holy Water
$ws.Range("H73").Value = 3705506.8
$ws.Range("I73").Value = 1034.3334
$ws.Range("J73").Value = 5557743
$ws.Range("K73").Value = 3103.0002
$ws.Range("L73").Value = 16673229
$ws.Range("M73").Value = -2167.0002
$ws.Range("N73").Value = -16675101

# Row 74: Adhesive of Antipathy | Wing Glue
$ws.Range("H74").Value = 2928.5
$ws.Range("I74").Value = 2721.25
$ws.Range("J74").Value = 3066.6667
$ws.Range("K74").Value = 2721.25
$ws.Range("L74").Value = 3066.6667
$ws.Range("M74").Value = -1785.25
$ws.Range("N74").Value = -4938.6667

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 2886.2554
$ws.Range("I76").Value = 2541.3333
$ws.Range("J76").Value = 3004.5144
$ws.Range("K76").Value = 2541.3333
$ws.Range("L76").Value = 3004.5144
$ws.Range("M76").Value = -2226.3333
$ws.Range("N76").Value = -3634.5144

# Row 77: It's Gonna Grow Back (L) | Wing Glue
$ws.Range("H77").Value = 2928.5
$ws.Range("I77").Value = 2721.25
$ws.Range("J77").Value = 3066.6667
$ws.Range("K77").Value = 13606.25
$ws.Range("L77").Value = 15333.3335
$ws.Range("M77").Value = -8926.25
$ws.Range("N77").Value = -24693.3335

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 2886.2554
$ws.Range("I79").Value = 2541.3333
$ws.Range("J79").Value = 3004.5144
$ws.Range("K79").Value = 2541.3333
$ws.Range("L79").Value = 3004.5144
$ws.Range("M79").Value = -1449.3333
$ws.Range("N79").Value = -5188.5144

# Row 106: Making Your Mark | Enchanted Palladium Ink
$ws.Range("H106").Value = 1194.3103
$ws.Range("I106").Value = 978
$ws.Range("J106").Value = 1675
$ws.Range("K106").Value = 978
$ws.Range("L106").Value = 1675
$ws.Range("M106").Value = -347
$ws.Range("N106").Value = -2937

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 1833816.8
$ws.Range("I116").Value = 2200103
$ws.Range("J116").Value = 2385.7144
$ws.Range("K116").Value = 2200103
$ws.Range("L116").Value = 2385.7144
$ws.Range("M116").Value = -2196661
$ws.Range("N116").Value = -9269.714400000001

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2104.0967
$ws.Range("I138").Value = 1568.4445
$ws.Range("J138").Value = 2845.7693
$ws.Range("K138").Value = 4705.333500000001
$ws.Range("L138").Value = 8537.3079
$ws.Range("M138").Value = 434.6664999999994
$ws.Range("N138").Value = -18817.3079

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 13705172
$ws.Range("I32").Value = 16668476
$ws.Range("J32").Value = 28385.54
$ws.Range("K32").Value = 16668476
$ws.Range("L32").Value = 28385.54
$ws.Range("M32").Value = -16668189
$ws.Range("N32").Value = -28959.54

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 2058.6487
$ws.Range("I45").Value = 1907.862
$ws.Range("J45").Value = 2605.25
$ws.Range("K45").Value = 1907.862
$ws.Range("L45").Value = 2605.25
$ws.Range("M45").Value = -1530.862
$ws.Range("N45").Value = -3359.25

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1612.7894
$ws.Range("I61").Value = 2026.619
$ws.Range("J61").Value = 1101.5883
$ws.Range("K61").Value = 2026.619
$ws.Range("L61").Value = 1101.5883
$ws.Range("M61").Value = -1814.619
$ws.Range("N61").Value = -1525.5883

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 1448.4615
$ws.Range("I102").Value = 1202.7273
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 1202.7273
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = 419.2727
$ws.Range("N102").Value = -6044

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1656.4286
$ws.Range("I110").Value = 1655.4546
$ws.Range("J110").Value = 1660
$ws.Range("K110").Value = 1655.4546
$ws.Range("L110").Value = 1660
$ws.Range("M110").Value = 389.5454
$ws.Range("N110").Value = -5750

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1612.7894
$ws.Range("I136").Value = 2026.619
$ws.Range("J136").Value = 1101.5883
$ws.Range("K136").Value = 6079.857
$ws.Range("L136").Value = 3304.7649
$ws.Range("M136").Value = -3529.857
$ws.Range("N136").Value = -8404.7649

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2620.07
$ws.Range("I105").Value = 1351.35
$ws.Range("J105").Value = 2937.25
$ws.Range("K105").Value = 1351.35
$ws.Range("L105").Value = 2937.25
$ws.Range("M105").Value = 395.6500000000001
$ws.Range("N105").Value = -6431.25

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 3626.2727
$ws.Range("I107").Value = 3048.625
$ws.Range("J107").Value = 5166.6665
$ws.Range("K107").Value = 3048.625
$ws.Range("L107").Value = 5166.6665
$ws.Range("M107").Value = -1128.625
$ws.Range("N107").Value = -9006.666499999999

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2575.1667
$ws.Range("I134").Value = 2672.9092
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 8018.7276
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -5483.7276
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
# Row 42: Live Freelance or Die | Heavy Steel Lance
$ws.Range("H42").Value = 7096
$ws.Range("I42").Value = 5826.6665
$ws.Range("K42").Value = 5826.6665
$ws.Range("M42").Value = -5233.6665

# Row 44: Stay on Target | Yarzonshell Harpoon
$ws.Range("H44").Value = 32000
$ws.Range("J44").Value = 32000
$ws.Range("L44").Value = 32000
$ws.Range("N44").Value = -32884

# Row 56: Trident and Error | Cobalt Trident
$ws.Range("H56").Value = 8198.25
$ws.Range("I56").Value = 7597.6665
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 7597.6665
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = -6752.6665
$ws.Range("N56").Value = -11690

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 372101.6
$ws.Range("I105").Value = 527754.2
$ws.Range("J105").Value = 2426.625
$ws.Range("K105").Value = 527754.2
$ws.Range("L105").Value = 2426.625
$ws.Range("M105").Value = -526007.2
$ws.Range("N105").Value = -5920.625

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 1711.6333
$ws.Range("I132").Value = 959.2917
$ws.Range("J132").Value = 4721
$ws.Range("K132").Value = 2877.8751
$ws.Range("L132").Value = 14163
$ws.Range("M132").Value = -347.8751000000002
$ws.Range("N132").Value = -19223

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2398.2856
$ws.Range("I134").Value = 893.34784
$ws.Range("J134").Value = 5282.75
$ws.Range("K134").Value = 2680.04352
$ws.Range("L134").Value = 15848.25
$ws.Range("M134").Value = -145.0435200000002
$ws.Range("N134").Value = -20918.25

# Row 139: Weaving a Path | Acacia Spinning Wheel
$ws.Range("H139").Value = 38995
$ws.Range("J139").Value = 38995
$ws.Range("L139").Value = 38995
$ws.Range("N139").Value = -49275

$ws = $wb.Worksheets.Item("CUL")
# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 734104.75
$ws.Range("I113").Value = 1815228.2
$ws.Range("J113").Value = 485.2143
$ws.Range("K113").Value = 5445684.6
$ws.Range("L113").Value = 1455.6429
$ws.Range("M113").Value = -5443514.6
$ws.Range("N113").Value = -5795.6429

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 952.1449
$ws.Range("J131").Value = 1003.93445
$ws.Range("L131").Value = 3011.80335
$ws.Range("N131").Value = -13091.80335

$ws = $wb.Worksheets.Item("GSM")
# Row 15: The Tusk at Hand | Fang Earrings
$ws.Range("H15").Value = 17000
$ws.Range("J15").Value = 17000
$ws.Range("L15").Value = 17000
$ws.Range("N15").Value = -17576

# Row 81: The Grander Temple | Dragon Fang Earrings
$ws.Range("H81").Value = 17000
$ws.Range("J81").Value = 17000
$ws.Range("L81").Value = 17000
$ws.Range("N81").Value = -18996

# Row 84: Man with a Dragon Earring (L) | Dragon Fang Earrings
$ws.Range("H84").Value = 17000
$ws.Range("J84").Value = 17000
$ws.Range("L84").Value = 51000
$ws.Range("N84").Value = -60984

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 770448.4
$ws.Range("I97").Value = 910410.8
$ws.Range("J97").Value = 655
$ws.Range("K97").Value = 910410.8
$ws.Range("L97").Value = 655
$ws.Range("M97").Value = -909914.8
$ws.Range("N97").Value = -1647

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 1847.0741
$ws.Range("I102").Value = 1884.3334
$ws.Range("J102").Value = 1716.6666
$ws.Range("K102").Value = 1884.3334
$ws.Range("L102").Value = 1716.6666
$ws.Range("M102").Value = -262.3334
$ws.Range("N102").Value = -4960.6666

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 6865.391
$ws.Range("I113").Value = 1732.2858
$ws.Range("J113").Value = 14850.223
$ws.Range("K113").Value = 1732.2858
$ws.Range("L113").Value = 14850.223
$ws.Range("M113").Value = 437.7141999999999
$ws.Range("N113").Value = -19190.223

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 3479.5293
$ws.Range("I122").Value = 2609.5
$ws.Range("J122").Value = 4722.4287
$ws.Range("K122").Value = 7828.5
$ws.Range("L122").Value = 14167.2861
$ws.Range("M122").Value = -5378.5
$ws.Range("N122").Value = -19067.2861

$ws = $wb.Worksheets.Item("LTW")
# Row 34: Breeches Served Cold | Goatskin Breeches
$ws.Range("H34").Value = 4281.5
$ws.Range("I34").Value = 3300
$ws.Range("J34").Value = 5263
$ws.Range("K34").Value = 3300
$ws.Range("L34").Value = 5263
$ws.Range("M34").Value = -3128
$ws.Range("N34").Value = -5607

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 1251750.2
$ws.Range("I68").Value = 1430286
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1430286
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1429537
$ws.Range("N68").Value = -3498

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 1251750.2
$ws.Range("I71").Value = 1430286
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 7151430
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -7147686
$ws.Range("N71").Value = -17488

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 1400.1786
$ws.Range("I136").Value = 1227.2727
$ws.Range("J136").Value = 2034.1666
$ws.Range("K136").Value = 3681.8181
$ws.Range("L136").Value = 6102.4998
$ws.Range("M136").Value = -1131.8181
$ws.Range("N136").Value = -11202.4998

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 20350.814
$ws.Range("I136").Value = 23068.348
$ws.Range("J136").Value = 4725
$ws.Range("K136").Value = 69205.04400000001
$ws.Range("L136").Value = 14175
$ws.Range("M136").Value = -19275
